# Insert a new weekly record for "Femacal de La Calera - Acelga" at row 274.
# Excel shifts all the existing rows (old 274..292) down to (275..293),
# extending the used range from A1:R292 to A1:R293.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above current row 274 (pushes 274..292 down to 275..293)
$ws.Rows.Item(274).Insert()

# Populate the new row 274 with the new weekly observation
$ws.Range("A274").Value = 3
$ws.Range("B274").Value = "Femacal de La Calera"
$ws.Range("C274").Value = "Coquimbo"
$ws.Range("D274").Value = 44610
$ws.Range("E274").Value = 5
$ws.Range("F274").Value = 100112009
$ws.Range("G274").Value = "Acelga"
$ws.Range("H274").Value = "Sin especificar"
$ws.Range("I274").Value = "Primera"
$ws.Range("J274").Value = 230
$ws.Range("K274").Value = 2500
$ws.Range("L274").Value = 2800
$ws.Range("M274").Value = 2643
$ws.Range("N274").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O274").Value = "Provincia de Quillota"
$ws.Range("P274").Value = 440
$ws.Range("Q274").Value = 6
$ws.Range("R274").Value = "Hortaliza"
